$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BE1").Value = 0.66281595854511
$ws.Range("BK2").Value = 0.96493453569194698
$ws.Range("BP2").Value = 0.759969748191482
$ws.Range("F4").Value = 0.74332417222778224
$ws.Range("U4").Value = 0.61805998820022645
$ws.Range("Q5").Value = 0.94674508219734466
$ws.Range("BD5").Value = 0.68126087412089043
$ws.Range("BF5").Value = 0.83982828501328533
$ws.Range("BK5").Value = 0.90175408978815796
$ws.Range("U6").Value = 0.75472400271803242
$ws.Range("AF6").Value = 0.87922494330213952
$ws.Range("G9").Value = 0.89359195201039676
$ws.Range("Z9").Value = 0.74343287628706256
$ws.Range("AA9").Value = 0.97379899406761816
$ws.Range("K10").Value = 0.99963404493445185
$ws.Range("AG10").Value = 0.70149553643529416
$ws.Range("BM11").Value = 0.97792002031481973
$ws.Range("BN11").Value = 0.96634272826445922
$ws.Range("D12").Value = 0.65518477103860451
$ws.Range("AE12").Value = 0.77776951892419754
$ws.Range("AP12").Value = 0.96986740233655189
$ws.Range("AX12").Value = 0.68226150584156198
$ws.Range("AW13").Value = 0.89972088699944663
$ws.Range("D15").Value = 0.98876918415699566
$ws.Range("M15").Value = 0.83064659876341129
$ws.Range("AD15").Value = 0.97088870389975845
$ws.Range("G16").Value = 0.81755918788411119
$ws.Range("Y16").Value = 0.98374777364116817
$ws.Range("BM16").Value = 0.72115573146959433
$ws.Range("S17").Value = 0.85769361861842408
$ws.Range("AT17").Value = 0.95083515541104102
$ws.Range("AC19").Value = 0.96420454327134619
$ws.Range("S20").Value = 0.63250669595132969
$ws.Range("V20").Value = 0.88477388592362605
$ws.Range("J21").Value = 0.64633864943796038
$ws.Range("L21").Value = 0.90791688908958235
$ws.Range("N21").Value = 0.85877850391074695
$ws.Range("T21").Value = 0.73595128299014934
$ws.Range("AT21").Value = 0.97810096343042385
$ws.Range("W22").Value = 0.98286666999799865
$ws.Range("AD22").Value = 0.9114335807857733
$ws.Range("BB22").Value = 0.83487906315406324
$ws.Range("X23").Value = 0.74737458944996082
$ws.Range("R24").Value = 0.96893932239349079
$ws.Range("BG24").Value = 0.71345579535571868
$ws.Range("K26").Value = 0.95698714821638187
$ws.Range("X26").Value = 0.93817839443504325
$ws.Range("AC27").Value = 0.83418729747969089
$ws.Range("AA28").Value = 0.82582123249017936
$ws.Range("Z29").Value = 0.96027134859073549
$ws.Range("F30").Value = 0.6864514300851996
$ws.Range("AA30").Value = 0.8192266641628152
$ws.Range("AD31").Value = 0.88041433195256924
$ws.Range("Z32").Value = 0.90146335895918406
$ws.Range("G33").Value = 0.88913603797515806
$ws.Range("I34").Value = 0.91364110686997624
$ws.Range("R34").Value = 0.92148917512100947
$ws.Range("AK34").Value = 0.82412114925553337
$ws.Range("P35").Value = 0.97657852131324807
$ws.Range("W35").Value = 0.83955867408072804
$ws.Range("AK35").Value = 0.92759347232146361
$ws.Range("BG35").Value = 0.76180967137457833
$ws.Range("AK36").Value = 0.9971358011446747
$ws.Range("V37").Value = 0.85440799684556845
$ws.Range("L38").Value = 0.73568687427423018
$ws.Range("AR38").Value = 0.97821066709211379
$ws.Range("A39").Value = 0.77204701879322601
$ws.Range("Q39").Value = 0.85012154853375921
$ws.Range("AE39").Value = 0.906008724138055
$ws.Range("AU39").Value = 0.96046063137089055
$ws.Range("C40").Value = 0.73590540600518961
$ws.Range("AP40").Value = 0.98096925509747224
$ws.Range("AA41").Value = 0.87296214794959792
$ws.Range("AN41").Value = 0.89065536383570709
$ws.Range("H42").Value = 0.61567134584442496
$ws.Range("AB42").Value = 0.70070723270434565
$ws.Range("AF42").Value = 0.77255822965843479
$ws.Range("I43").Value = 0.86319250866528541
$ws.Range("U43").Value = 0.9014500478087738
$ws.Range("AB43").Value = 0.80968608466024716
$ws.Range("H44").Value = 0.93746267050973198
$ws.Range("J44").Value = 0.75793903886626146
$ws.Range("AB45").Value = 0.84682661959966077
$ws.Range("AJ45").Value = 0.77280073216705225
$ws.Range("AK45").Value = 0.90278474575587686
$ws.Range("AH46").Value = 0.59889648604301216
$ws.Range("AL46").Value = 0.93598950434121897
$ws.Range("AW46").Value = 0.80419017149395156
$ws.Range("AF47").Value = 0.7739687553748974
$ws.Range("AQ48").Value = 0.69537254506232782
$ws.Range("N49").Value = 0.9741061276899432
$ws.Range("AR49").Value = 0.71549926925853158
$ws.Range("BF49").Value = 0.72666537229094708
$ws.Range("AV50").Value = 0.99518315528541312
$ws.Range("H51").Value = 0.86971259556961034
$ws.Range("Y51").Value = 0.87591123820896799
$ws.Range("AQ51").Value = 0.79377681740287054
$ws.Range("AX51").Value = 0.53736805185336434
$ws.Range("BJ51").Value = 0.74012751188476067
$ws.Range("BO51").Value = 0.91493552171207515
$ws.Range("AC52").Value = 0.94147617986203025
$ws.Range("BB52").Value = 0.90736234154555961
$ws.Range("BG52").Value = 0.65728406790301774
$ws.Range("A53").Value = 0.94825683812880301
$ws.Range("BB53").Value = 0.65410475353242337
$ws.Range("BL53").Value = 0.99144926987561599
$ws.Range("G55").Value = 0.8637685516924456
$ws.Range("AM55").Value = 0.8657474703057555
$ws.Range("BP55").Value = 0.98797777745575588
$ws.Range("H56").Value = 0.87962730787233112
$ws.Range("AU56").Value = 0.94956728574409244
$ws.Range("AW56").Value = 0.92576107359726856
$ws.Range("B57").Value = 0.92796070053328772
$ws.Range("W57").Value = 0.81928921049911119
$ws.Range("AX57").Value = 0.65481945478929904
$ws.Range("F58").Value = 0.81261525351962327
$ws.Range("C59").Value = 0.82205794333741156
$ws.Range("AT59").Value = 0.93763796931202059
$ws.Range("U60").Value = 0.87764944312677784
$ws.Range("X60").Value = 0.72031193255310377
$ws.Range("Y60").Value = 0.95061406965925421
$ws.Range("AS60").Value = 0.98571833090996974
$ws.Range("AU60").Value = 0.84599952585384963
$ws.Range("R61").Value = 0.93778995220522232
$ws.Range("AV61").Value = 0.95440502748421574
$ws.Range("BD61").Value = 0.85090475066042626
$ws.Range("M62").Value = 0.60634150291596201
$ws.Range("BB62").Value = 0.79076835124792866
$ws.Range("AL63").Value = 0.79563803679004963
$ws.Range("BC63").Value = 0.77311626856949922
$ws.Range("B64").Value = 0.61831048933778876
$ws.Range("BO64").Value = 0.98696517829904207
$ws.Range("BK65").Value = 0.84349530374311832
$ws.Range("AF66").Value = 0.96973425659115986
$ws.Range("AO66").Value = 0.79441931088027851
$ws.Range("BO68").Value = 0.74468556446497569
